$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.277.21'
$ws.Range("E2").Value = '  -0.12%  '

$ws.Range("D3").Value = '1.929.53'

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.002'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.12%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '248.88'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.17%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.7143'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.23%  '

$ws.Range("E7").Value = '  +0.17%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3214'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -1.83%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '27.37'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.41%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07098'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +4.10%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.7922'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -1.56%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.08037'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.18%  '

$ws.Range("D13").Value = '1.926.55'
$ws.Range("E13").Value = '  -0.29%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.376'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.63%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '94.59'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.24%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '14.68'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +1.37%  '

$ws.Range("D17").Value = '30.271.71'

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '256.23'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +1.10%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.000008029'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.54%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '5.770'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.66%  '

$ws.Range("D21").Value = '2.181.43'
$ws.Range("E21").Value = '  -0.28%  '

$ws.Range("E22").Value = '  +0.19%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.001'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.11%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.816'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.87%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.562'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.95%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '166.08'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +4.42%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '19.11'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.10%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.279'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -4.80%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.1278'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -3.66%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.369'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +2.40%  '

$ws.Range("E31").Value = '  -1.54%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.394'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.06%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.134'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -1.04%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.05174'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +2.32%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.262'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +3.59%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7457'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +1.11%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.769'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.60%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01957'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.56%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.808'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.68%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '77.72'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -1.85%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '6.340'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -3.89%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.4486'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.78%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.981'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.02%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.8458'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +1.44%  '

$ws.Range("E45").Value = '  +0.11%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '100.78'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -1.21%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '9.677'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -1.05%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '7.428'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +2.05%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '36.40'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.15%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.06117'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +2.97%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.4165'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +2.61%  '

